# Atualizacao de bases das ligas, do dia: 17-02-2024 as 22:47
# Australia ALeague - apply row updates per upstream diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 97 and 98: the two matches' data (every column except the running
#    index in column A) were swapped between the rows.
#    Row 97 (A97=95) now carries what used to be row 98's data (id 5404732,
#    Melbourne City vs Western Sydney Wanderers, ...).
#    Row 98 (A98=96) now carries what used to be row 97's data (id 5400063,
#    Central Coast Mariners vs Adelaide United, ...).
# ---------------------------------------------------------------------------

# New row 97 content (previously row 98's values)
$ws.Range("B97").Value = 5404732
$ws.Range("F97").Value = "Adelaide United"
$ws.Range("G97").Value = "Central Coast Mariners"
$ws.Range("H97").Value = 1
$ws.Range("I97").Value = 4
$ws.Range("J97").Value = "A"
$ws.Range("K97").Value = 2.3
$ws.Range("L97").Value = 3.75
$ws.Range("M97").Value = 2.75
$ws.Range("N97").Value = 2.625
$ws.Range("O97").Value = 4
$ws.Range("P97").Value = 2.4
$ws.Range("Q97").Value = 0
$ws.Range("R97").Value = 2.05
$ws.Range("S97").Value = 1.8
$ws.Range("T97").Value = 3.75
$ws.Range("U97").Value = 2
$ws.Range("V97").Value = 1.85
$ws.Range("W97").Value = -1
$ws.Range("X97").Value = -1
$ws.Range("Y97").Value = 1.4
$ws.Range("Z97").Value = -1
$ws.Range("AA97").Value = 0.8
$ws.Range("AB97").Value = 1
$ws.Range("AC97").Value = -1

# New row 98 content (previously row 97's values)
$ws.Range("B98").Value = 5400063
$ws.Range("F98").Value = "Melbourne City"
$ws.Range("G98").Value = "Western Sydney Wanderers"
$ws.Range("H98").Value = 3
$ws.Range("I98").Value = 2
$ws.Range("J98").Value = "H"
$ws.Range("K98").Value = 1.75
$ws.Range("L98").Value = 3.8
$ws.Range("M98").Value = 4
$ws.Range("N98").Value = 2
$ws.Range("O98").Value = 4
$ws.Range("P98").Value = 3.4
$ws.Range("Q98").Value = -0.5
$ws.Range("R98").Value = 2.025
$ws.Range("S98").Value = 1.825
$ws.Range("T98").Value = 3
$ws.Range("U98").Value = 1.85
$ws.Range("V98").Value = 2
$ws.Range("W98").Value = 1
$ws.Range("X98").Value = -1
$ws.Range("Y98").Value = -1
$ws.Range("Z98").Value = 1.025
$ws.Range("AA98").Value = -1
$ws.Range("AB98").Value = 0.8500000000000001
$ws.Range("AC98").Value = -1

# ---------------------------------------------------------------------------
# 2) Row 208: refreshed closing Asian-handicap odds.
# ---------------------------------------------------------------------------
$ws.Range("R208").Value = 2.02
$ws.Range("S208").Value = 1.88
$ws.Range("U208").Value = 1.85
$ws.Range("V208").Value = 2

# ---------------------------------------------------------------------------
# 3) Row 209: refreshed closing Asian-handicap odds.
# ---------------------------------------------------------------------------
$ws.Range("R209").Value = 1.87
$ws.Range("S209").Value = 2.03
$ws.Range("U209").Value = 1.95
$ws.Range("V209").Value = 1.9

# ---------------------------------------------------------------------------
# 4) New row 210: brand new fixture appended at the bottom of the table.
#    Copy formatting from the row above first, so the running-index cell
#    (A210) and the date cell (E210) pick up the same styles as the rest of
#    the sheet (bold/centered/bordered index column, custom date format).
# ---------------------------------------------------------------------------
$ws.Range("A209").Copy()
$ws.Range("A210").PasteSpecial(-4122)
$ws.Range("E209").Copy()
$ws.Range("E210").PasteSpecial(-4122)

$ws.Range("A210").Value = 208
$ws.Range("B210").Value = 7661946
$ws.Range("C210").Value = "Australia ALeague"
$ws.Range("D210").Value = "Australia ALeague"
$ws.Range("E210").Value = 45342.20833333334
$ws.Range("F210").Value = "Melbourne Victory"
$ws.Range("G210").Value = "Western United FC"
$ws.Range("K210").Value = 1.45
$ws.Range("L210").Value = 4.75
$ws.Range("M210").Value = 6.5
$ws.Range("N210").Value = 1.45
$ws.Range("O210").Value = 4.75
$ws.Range("P210").Value = 6.5
$ws.Range("Q210").Value = -1.25
$ws.Range("R210").Value = 2
$ws.Range("S210").Value = 1.9
$ws.Range("T210").Value = 3
$ws.Range("U210").Value = 1.875
$ws.Range("V210").Value = 1.975
$ws.Range("W210").Value = 0
$ws.Range("X210").Value = 0
$ws.Range("Y210").Value = 0
$ws.Range("Z210").Value = 0
$ws.Range("AA210").Value = 0

Write-Host "edit complete"
